# Add data for 2022-11-24
# Updates violent-crime YTD counts across the Citywide Totals rollup, the
# By Neighborhood rollup, and the individual neighborhood sheets that saw
# new/recategorized incidents.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("C2").Value = 69
$ws.Range("E2").Value = 69
$ws.Range("I2").Value = 113
$ws.Range("H3").Value = 150
$ws.Range("B6").Value = 364
$ws.Range("D6").Value = 396
$ws.Range("E6").Value = 450
$ws.Range("F6").Value = 501
$ws.Range("G6").Value = 428
$ws.Range("H6").Value = 432
$ws.Range("B7").Value = 488
$ws.Range("C7").Value = 612
$ws.Range("D7").Value = 622
$ws.Range("E7").Value = 671
$ws.Range("F7").Value = 727
$ws.Range("G7").Value = 655
$ws.Range("H7").Value = 699
$ws.Range("I7").Value = 810

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("B7").Value = 19
$ws.Range("F8").Value = 47
$ws.Range("E29").Value = 8
$ws.Range("E32").Value = 62
$ws.Range("E35").Value = 6
$ws.Range("D36").Value = 34
$ws.Range("E36").Value = 36
$ws.Range("G36").Value = 26
$ws.Range("I43").Value = 6
$ws.Range("G50").Value = 15
$ws.Range("B63").Value = 7
$ws.Range("C65").Value = 21
$ws.Range("D65").Value = 23
$ws.Range("H69").Value = 3
$ws.Range("H74").Value = 15
$ws.Range("B98").Value = 488
$ws.Range("C98").Value = 612
$ws.Range("D98").Value = 622
$ws.Range("E98").Value = 671
$ws.Range("F98").Value = 727
$ws.Range("G98").Value = 655
$ws.Range("H98").Value = 699
$ws.Range("I98").Value = 810

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("B5").Value = 17
$ws.Range("B6").Value = 19
$ws.Range("I2").Value = 2

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I6").Value = 6

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("F5").Value = 33
$ws.Range("F6").Value = 47

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("E6").Value = 49
$ws.Range("E7").Value = 62

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("E2").Value = 5
$ws.Range("D6").Value = 19
$ws.Range("G6").Value = 14
$ws.Range("D7").Value = 34
$ws.Range("E7").Value = 36
$ws.Range("G7").Value = 26

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("G5").Value = 14
$ws.Range("G6").Value = 15

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("E5").Value = 5
$ws.Range("E6").Value = 6

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("C2").Value = 3
$ws.Range("D5").Value = 22
$ws.Range("C6").Value = 21
$ws.Range("D6").Value = 23

$ws = $wb.Worksheets.Item("River North")
$ws.Range("H5").Value = 11
$ws.Range("H6").Value = 15

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("E5").Value = 4
$ws.Range("E6").Value = 8

$ws = $wb.Worksheets.Item("New City")
$ws.Range("B4").Value = 5
$ws.Range("B5").Value = 7

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("E3").Value = 2
$ws.Range("E5").Value = 3
